$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.03955433333333334
$ws.Range("H2").Value = 0.118663
$ws.Range("M2").Value = 361.2779286666667
$ws.Range("N2").Value = 1083.833786
$ws.Range("O2").Value = 0.9679392703861037
$ws.Range("P2").Value = 0.9679392703861038
$ws.Range("Q2").Value = 14.29010761645756
$ws.Range("R2").Value = 128.610968548118
$ws.Range("S2").Value = 0.9679392703861037
$ws.Range("T2").Value = 0.9679392703861038

# Row 3
$ws.Range("G3").Value = 0.03955433333333334
$ws.Range("H3").Value = 0.118663
$ws.Range("O3").Value = 0.015995373883918
$ws.Range("P3").Value = 0.015995373883918
$ws.Range("Q3").Value = 0.2361466479973333
$ws.Range("R3").Value = 2.125319831976
$ws.Range("S3").Value = 0.015995373883918
$ws.Range("T3").Value = 0.015995373883918

# Row 4
$ws.Range("G4").Value = 0.03955433333333334
$ws.Range("H4").Value = 0.118663
$ws.Range("M4").Value = 3.890485666666667
$ws.Range("N4").Value = 11.671457
$ws.Range("O4").Value = 0.01042342628440887
$ws.Range("P4").Value = 0.01042342628440887
$ws.Range("Q4").Value = 0.1538855668878889
$ws.Range("R4").Value = 1.384970101991
$ws.Range("S4").Value = 0.01042342628440887
$ws.Range("T4").Value = 0.01042342628440887

# Row 5
$ws.Range("G5").Value = 0.03955433333333334
$ws.Range("H5").Value = 0.118663
$ws.Range("M5").Value = 2.105818666666667
$ws.Range("N5").Value = 6.317456
$ws.Range("O5").Value = 0.005641929445569353
$ws.Range("P5").Value = 0.005641929445569354
$ws.Range("Q5").Value = 0.0832942534808889
$ws.Range("R5").Value = 0.749648281328
$ws.Range("S5").Value = 0.005641929445569353
$ws.Range("T5").Value = 0.005641929445569354
